$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change: D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen), S (Precio $/Kg)
# The edit rotates these values across rows 2,3,5,6,7 (row 4 untouched):
#   new row2 <- old row6
#   new row3 <- old row7
#   new row5 <- old row3
#   new row6 <- old row2
#   new row7 <- old row5

function Get-RowData($row) {
    return @{
        D = $ws.Cells.Item($row, 4).Value2
        L = $ws.Cells.Item($row, 12).Value2
        M = $ws.Cells.Item($row, 13).Value2
        N = $ws.Cells.Item($row, 14).Value2
        O = $ws.Cells.Item($row, 15).Value2
        P = $ws.Cells.Item($row, 16).Value2
        R = $ws.Cells.Item($row, 18).Value2
        S = $ws.Cells.Item($row, 19).Value2
    }
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 4).Value2 = $data.D
    $ws.Cells.Item($row, 12).Value2 = $data.L
    $ws.Cells.Item($row, 13).Value2 = $data.M
    $ws.Cells.Item($row, 14).Value2 = $data.N
    $ws.Cells.Item($row, 15).Value2 = $data.O
    $ws.Cells.Item($row, 16).Value2 = $data.P
    $ws.Cells.Item($row, 18).Value2 = $data.R
    $ws.Cells.Item($row, 19).Value2 = $data.S
}

$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row5 = Get-RowData 5
$row6 = Get-RowData 6
$row7 = Get-RowData 7

Set-RowData 2 $row6
Set-RowData 3 $row7
Set-RowData 5 $row3
Set-RowData 6 $row2
Set-RowData 7 $row5
